$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.970.73"
$ws.Range("E2").Value = "  -3.31%  "
$ws.Range("D3").Value = "1.794.92"
$ws.Range("E3").Value = "  -3.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4193"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3561"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07083"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8441"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("E11").Value = "  -5.28%  "
$ws.Range("D12").Value = "1.803.96"
$ws.Range("E12").Value = "  -4.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.284"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.332"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06759"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008640"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").Value = "26.936.12"
$ws.Range("E21").Value = "  -3.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.053"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").Value = "2.005.91"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.931"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.64"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.07"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.93%  "
$ws.Range("E28").Value = "  -6.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.92"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.634"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -12.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08943"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7161"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -9.45%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.859"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.79%  "
$ws.Range("E34").Value = "  -7.57%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.072"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -9.38%  "
$ws.Range("E37").Value = "  -3.18%  "
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05095"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1623"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4935"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.568"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -9.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.966"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -12.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.013"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.42"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06296"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4499"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.590"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "61.82"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.34%  "
